$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace a paragraph's full range (text+runs) with a new WordML
# <w:p> fragment via Range.InsertXML, which accepts a raw OOXML fragment and
# splices it into the document in place of the given range.
# ---------------------------------------------------------------------------
function Set-ParagraphXml($paraIndex, $fragment) {
    $para = $d.Paragraphs($paraIndex)
    $rng = $para.Range
    $rng.InsertXML($fragment)
}

$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------------------
# 1) "Team Members:" -> "Pre" + "Team" (wrapped in spellStart/spellEnd) +
#    " Members:", i.e. the heading becomes "PreTeam Members: Project Group 1"
# ---------------------------------------------------------------------------
$p1 = @"
<w:p $wns>
  <w:pPr>
    <w:spacing w:line="240" w:lineRule="auto"/>
    <w:contextualSpacing/>
    <w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/></w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>Pre</w:t></w:r>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>Team</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> Members:</w:t></w:r>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> Project </w:t></w:r>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/></w:rPr><w:t>Group 1</w:t></w:r>
</w:p>
"@
Set-ParagraphXml 1 $p1

# ---------------------------------------------------------------------------
# 19) "Download data from Kaggle - Matt" -> wrap "Matt" in gramStart/gramEnd
# ---------------------------------------------------------------------------
$p19 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr>
    <w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t>Download</w:t></w:r>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t xml:space="preserve"> data </w:t></w:r>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t>from Kaggle</w:t></w:r>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t xml:space="preserve"> - </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t>Matt</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
</w:p>
"@
Set-ParagraphXml 19 $p19

# ---------------------------------------------------------------------------
# 21) "Use Python to clean and format dataset" -> wrap "dataset" in
#     gramStart/gramEnd
# ---------------------------------------------------------------------------
$p21 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr>
    <w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t xml:space="preserve">Use </w:t></w:r>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t>Python</w:t></w:r>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t xml:space="preserve"> to clean and format </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t>dataset</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
</w:p>
"@
Set-ParagraphXml 21 $p21

# ---------------------------------------------------------------------------
# 22) "Inspect dataset in Pandas data frame" -> wrap "Pandas" in
#     gramStart/gramEnd
# ---------------------------------------------------------------------------
$p22 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr>
    <w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t xml:space="preserve">Inspect dataset in </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t>Pandas</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t xml:space="preserve"> data frame</w:t></w:r>
</w:p>
"@
Set-ParagraphXml 22 $p22

# ---------------------------------------------------------------------------
# 23) "Machine Learning  - Matt, Luis, Emily" -> re-split into
#     "Machine " + gramStart + "Learning " + " -" + gramEnd + " " + "Matt, Luis, Emily"
# ---------------------------------------------------------------------------
$p23 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr>
    <w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t xml:space="preserve">Machine </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t xml:space="preserve">Learning </w:t></w:r>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t xml:space="preserve"> &#8211;</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t>Matt, Luis, Emily</w:t></w:r>
</w:p>
"@
Set-ParagraphXml 23 $p23

# ---------------------------------------------------------------------------
# 24) "Create & run code using Python" -> wrap "Python" in gramStart/gramEnd
# ---------------------------------------------------------------------------
$p24 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr>
    <w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t>Create &amp; run code</w:t></w:r>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t xml:space="preserve"> using </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t>Python</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
</w:p>
"@
Set-ParagraphXml 24 $p24

# ---------------------------------------------------------------------------
# 27) "Create a write-up summarizing major findings and implications. - Emily "
#     -> "... - All" (the trailing name changes from Emily to All)
# ---------------------------------------------------------------------------
$p27 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr>
    <w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t>Create a write-up summarizing major findings and implications.</w:t></w:r>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t xml:space="preserve"> &#8211; </w:t></w:r>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t>All</w:t></w:r>
</w:p>
"@
Set-ParagraphXml 27 $p27

# ---------------------------------------------------------------------------
# 29) "Instructions on how to use and interact with the project" -> wrap
#     "project" in gramStart/gramEnd
# ---------------------------------------------------------------------------
$p29 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr>
    <w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t xml:space="preserve">Instructions on how to use and interact with the </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t>project</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
</w:p>
"@
Set-ParagraphXml 29 $p29

# ---------------------------------------------------------------------------
# 31) "References for any code used that is not your own" -> wrap "own" in
#     gramStart/gramEnd
# ---------------------------------------------------------------------------
$p31 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr>
    <w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t xml:space="preserve">References for any code used that is not your </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t>own</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
</w:p>
"@
Set-ParagraphXml 31 $p31

Write-Output "Done."
